# Weekly fruit/vegetable price update: insert a new observation row at
# row 109 (pushing the existing rows 109:186 down to 110:187) and populate
# it with the new week's data. All other rows keep their original values,
# just shifted down by one row - which Rows.Insert() handles natively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 109; rows 109:186 shift down to 110:187.
$ws.Rows("109:109").Insert()

# Populate the newly inserted row 109 with the new record.
$ws.Cells.Item(109, 1).Value = 11
$ws.Cells.Item(109, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(109, 3).Value = "Bíobío"
$ws.Cells.Item(109, 4).Value = 44978
$ws.Cells.Item(109, 5).Value = 8
$ws.Cells.Item(109, 6).Value = 100112043
$ws.Cells.Item(109, 7).Value = "Pepino ensalada"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 100
$ws.Cells.Item(109, 11).Value = 7500
$ws.Cells.Item(109, 12).Value = 8000
$ws.Cells.Item(109, 13).Value = 7750
$ws.Cells.Item(109, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(109, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(109, 16).Value = 129
$ws.Cells.Item(109, 17).Value = 60
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D (style index 2,
# numFmtId 165 date/time) - Rows.Insert() already carries that format down
# from the row above, but set it explicitly to be safe.
$ws.Cells.Item(109, 4).NumberFormat = $ws.Cells.Item(110, 4).NumberFormat
